# Update the "想去人数" (number of people interested) counts that were
# regenerated for this gh-pages data refresh.

$wb = $excel.ActiveWorkbook

# Sheet "展览"
$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F4").Value = 7748
$wsExhibit.Range("F5").Value = 5643
$wsExhibit.Range("F11").Value = 240

# Sheet "演出"
$wsShow = $wb.Worksheets.Item("演出")
$wsShow.Range("F2").Value = 86

# Sheet "全部类型"
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F4").Value = 7748
$wsAll.Range("F5").Value = 5643
$wsAll.Range("F11").Value = 86
$wsAll.Range("F13").Value = 240
